# Adiciona o item "HTML, CSS e Java Script" a lista de Habilidades,
# logo apos "Versionamento com GitHub" e antes da linha horizontal
# que separa a secao de Habilidades da secao de Idiomas.

$d = $word.ActiveDocument

# Localiza o paragrafo que contem "Versionamento com GitHub" (ultimo
# item da lista de Habilidades) percorrendo os paragrafos do documento.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Versionamento com GitHub*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Nao foi possivel localizar o paragrafo 'Versionamento com GitHub'."
}

$anchorPara = $d.Paragraphs.Item($targetIndex)
$anchorRange = $anchorPara.Range

# Colapsa o range para o final do paragrafo-ancora e insere um novo
# paragrafo logo depois dele (antes da linha horizontal seguinte).
$anchorRange.Collapse(0)
$anchorRange.InsertParagraphAfter()

# O novo paragrafo criado herda a formatacao de lista numerada
# (numPr / numId=13) do paragrafo anterior automaticamente.
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "HTML, CSS e Java Script"
